$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.315672755241394
$ws.Range("B1").Value = 1.478198528289795
$ws.Range("C1").Value = 3.792777299880981
$ws.Range("D1").Value = 3.434493780136108
$ws.Range("E1").Value = 1.026684403419495
